$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert a new "2022-Q1" sheet right after "2021-Q2" (and before "总计").
#    Copying "2021-Q2" gives us an identical header/style template (bold
#    bordered header row, centered "A" index column) for free.
# ---------------------------------------------------------------------------
$src = $wb.Worksheets.Item("2021-Q2")
$src.Copy($null, $src)
$newSheet = $wb.Worksheets.Item(3)
$newSheet.Name = "2022-Q1"

# The template only carries one data row; duplicate it so we have a second
# row (row 3) with the same styling (keeps the "A" column's centered style).
$newSheet.Range("A2:H2").Copy($newSheet.Range("A3:H3"))

# Header row: only the "基金规模" column label changed from the template.
$newSheet.Range("D1").Value = "基金规模"

# --- Row 2: 001303 / 银华稳利灵活配置混合A ---------------------------------
$newSheet.Range("A2").Value = 0
$row2 = $newSheet.Range("B2:G2")
$row2.NumberFormat = "@"
$newSheet.Range("B2").Value = "001303"
$newSheet.Range("C2").Value = "银华稳利灵活配置混合A"
$newSheet.Range("D2").Value = "0.18"
$newSheet.Range("E2").Value = "28.88"
$newSheet.Range("F2").Value = "0.76"
$newSheet.Range("G2").Value = "0.0014"
$row2.ClearFormats()
$newSheet.Range("H2").Value = 3

# --- Row 3: 002323 / 银华稳利灵活配置混合C ---------------------------------
$newSheet.Range("A3").Value = 1
$row3 = $newSheet.Range("B3:G3")
$row3.NumberFormat = "@"
$newSheet.Range("B3").Value = "002323"
$newSheet.Range("C3").Value = "银华稳利灵活配置混合C"
$newSheet.Range("D3").Value = "0.12"
$newSheet.Range("E3").Value = "28.88"
$newSheet.Range("F3").Value = "0.76"
$newSheet.Range("G3").Value = "0.0009"
$row3.ClearFormats()
$newSheet.Range("H3").Value = 3

# ---------------------------------------------------------------------------
# 2) Prepend a "2022-Q1" row to the "总计" summary sheet, pushing the
#    existing "2021-Q2"/"2020-Q4" rows down by one.
# ---------------------------------------------------------------------------
$sum = $wb.Worksheets.Item("总计")

$sum.Range("A3:D3").Copy($sum.Range("A4:D4"))
$sum.Range("A2:D2").Copy($sum.Range("A3:D3"))

$sum.Range("A4").Value = 2
$sum.Range("A3").Value = 1
$sum.Range("A2").Value = 0
$sum.Range("B2").Value = "2022-Q1"
$sum.Range("C2").Value = 2
$sum.Range("D2").Value = 0

# Restore the originally active sheet/selection (unrelated to the data edit).
$wb.Worksheets.Item("2020-Q4").Activate()
